$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change C27's number format from date-only (YYYY-MM-DD) to date-time
# (YYYY-MM-DD HH:MM:SS), while keeping its existing value (45758).
$ws.Range("C27").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append a new data row 28 with values matching row 20-23's data point.
$ws.Range("A28").Value = 781.86
$ws.Range("B28").Value = 679.38
$ws.Range("C28").Value = 45754
$ws.Range("C28").NumberFormat = "YYYY-MM-DD"
